$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '67.812.55'
Set-TextValue $ws.Range("E2") '  +0.25%  '

Set-TextValue $ws.Range("D3") '3.824.18'
Set-TextValue $ws.Range("E3") '  +1.12%  '

Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  -0.08%  '

Set-TextValue $ws.Range("D5") '604.87'
Set-TextValue $ws.Range("E5") '  +1.60%  '

Set-TextValue $ws.Range("D6") '167.14'
Set-TextValue $ws.Range("E6") '  +0.41%  '

Set-TextValue $ws.Range("E7") '  -0.13%  '

Set-TextValue $ws.Range("E8") '  +0.15%  '

Set-TextValue $ws.Range("E9") '  +0.67%  '

Set-TextValue $ws.Range("E10") '  -0.92%  '

Set-TextValue $ws.Range("E12") '  -0.58%  '

Set-TextValue $ws.Range("D13") '36.03'
Set-TextValue $ws.Range("E13") '  -0.71%  '

Set-TextValue $ws.Range("D14") '4.463.20'
Set-TextValue $ws.Range("E14") '  +1.02%  '

Set-TextValue $ws.Range("D15") '3.779.42'
Set-TextValue $ws.Range("E15") '  +0.03%  '

Set-TextValue $ws.Range("E16") '  +0.17%  '

Set-TextValue $ws.Range("D17") '67.827.00'
Set-TextValue $ws.Range("E17") '  +0.30%  '

Set-TextValue $ws.Range("D18") '7.09'
Set-TextValue $ws.Range("E18") '  +1.37%  '

Set-TextValue $ws.Range("E19") '  +0.53%  '

Set-TextValue $ws.Range("D20") '465.27'
Set-TextValue $ws.Range("E20") '  +1.75%  '

Set-TextValue $ws.Range("D21") '9.94'
Set-TextValue $ws.Range("E21") '  -1.70%  '

Set-TextValue $ws.Range("E22") '  +0.48%  '

Set-TextValue $ws.Range("D23") '0.0000150'
Set-TextValue $ws.Range("E23") '  -3.67%  '

Set-TextValue $ws.Range("D24") '83.45'
Set-TextValue $ws.Range("E24") '  +0.02%  '

Set-TextValue $ws.Range("D25") '12.08'
Set-TextValue $ws.Range("E25") '  +1.23%  '

Set-TextValue $ws.Range("E27") '  +0.19%  '

Set-TextValue $ws.Range("E28") '  +0.22%  '

Set-TextValue $ws.Range("D29") '3.970.46'
Set-TextValue $ws.Range("E29") '  +1.00%  '

Set-TextValue $ws.Range("E30") '  +0.33%  '

Set-TextValue $ws.Range("E31") '  +1.58%  '

Set-TextValue $ws.Range("E32") '  +1.90%  '

Set-TextValue $ws.Range("D33") '29.75'
Set-TextValue $ws.Range("E33") '  -0.29%  '

Set-TextValue $ws.Range("D34") '1.00'
Set-TextValue $ws.Range("E34") '  +0.05%  '

Set-TextValue $ws.Range("D35") '9.10'
Set-TextValue $ws.Range("E35") '  -1.35%  '

Set-TextValue $ws.Range("E36") '  -0.02%  '

Set-TextValue $ws.Range("D37") '3.36'
Set-TextValue $ws.Range("E37") '  +0.35%  '

Set-TextValue $ws.Range("E38") '  +0.11%  '

Set-TextValue $ws.Range("E39") '  +0.69%  '

Set-TextValue $ws.Range("E40") '  +0.73%  '

Set-TextValue $ws.Range("D41") '1.00'
Set-TextValue $ws.Range("E41") '  +0.01%  '

Set-TextValue $ws.Range("E42") '  +0.02%  '

Set-TextValue $ws.Range("D43") '48.11'
Set-TextValue $ws.Range("E43") '  +2.08%  '

Set-TextValue $ws.Range("E44") '  +0.54%  '

Set-TextValue $ws.Range("B45") 'Arweave'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws.Range("D45") '43.32'
Set-TextValue $ws.Range("E45") '  -4.58%  '

Set-TextValue $ws.Range("B46") 'EnergySwap'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("E46") '  +10.91%  '

Set-TextValue $ws.Range("D47") '1.41'
Set-TextValue $ws.Range("E47") '  +12.60%  '

Set-TextValue $ws.Range("D48") '8.36'
Set-TextValue $ws.Range("E48") '  +0.18%  '

Set-TextValue $ws.Range("D49") '148.04'

Set-TextValue $ws.Range("D50") '1.85'
Set-TextValue $ws.Range("E50") '  +0.40%  '

Set-TextValue $ws.Range("D51") '388.64'
Set-TextValue $ws.Range("E51") '  -0.33%  '

